$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a cell to hold a literal text value (matches the source file's
# inlineStr cells) even when the text looks numeric, without leaving a
# lingering custom style on the cell.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "58.435.68"
$ws.Range("E2").Value = "  +0.77%  "
Set-TextValue $ws.Range("D3") "2.488.10"
$ws.Range("E3").Value = "  +0.80%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue $ws.Range("D5") "519.75"
$ws.Range("E5").Value = "  +0.22%  "
Set-TextValue $ws.Range("D6") "132.01"
$ws.Range("E6").Value = "  +0.80%  "
Set-TextValue $ws.Range("D7") "0.997"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  -0.19%  "
Set-TextValue $ws.Range("D9") "2.522.59"
Set-TextValue $ws.Range("D10") "0.0976"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  -2.10%  "
Set-TextValue $ws.Range("D13") "0.333"
$ws.Range("E13").Value = "  -2.72%  "
Set-TextValue $ws.Range("D14") "2.936.17"
$ws.Range("E14").Value = "  +1.06%  "
Set-TextValue $ws.Range("D15") "58.255.50"
$ws.Range("E15").Value = "  +0.63%  "
Set-TextValue $ws.Range("D16") "22.23"
$ws.Range("E16").Value = "  -0.24%  "
Set-TextValue $ws.Range("D17") "0.0000135"
$ws.Range("E17").Value = "  -0.37%  "
Set-TextValue $ws.Range("D18") "2.508.76"
$ws.Range("E18").Value = "  +1.32%  "
Set-TextValue $ws.Range("D19") "10.73"
$ws.Range("E19").Value = "  -0.75%  "
Set-TextValue $ws.Range("D20") "323.48"
$ws.Range("E20").Value = "  +1.20%  "
Set-TextValue $ws.Range("D21") "4.17"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("E22").Value = "  +6.12%  "
Set-TextValue $ws.Range("D23") "0.996"
$ws.Range("E23").Value = "  -0.39%  "
Set-TextValue $ws.Range("D24") "63.58"
$ws.Range("E24").Value = "  -0.79%  "
Set-TextValue $ws.Range("D25") "0.406"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("E26").Value = "  +1.53%  "
Set-TextValue $ws.Range("D27") "0.992"
$ws.Range("E27").Value = "  -0.76%  "
Set-TextValue $ws.Range("D28") "7.38"
$ws.Range("E28").Value = "  +0.78%  "
Set-TextValue $ws.Range("D29") "0.0₃0748"
$ws.Range("E29").Value = "  -0.39%  "
Set-TextValue $ws.Range("D30") "168.99"
$ws.Range("E30").Value = "  +1.59%  "
Set-TextValue $ws.Range("D31") "1.70"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("E32").Value = "  +3.61%  "
Set-TextValue $ws.Range("D33") "6.28"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  -0.14%  "
Set-TextValue $ws.Range("D35") "0.995"
$ws.Range("E35").Value = "  -0.30%  "
Set-TextValue $ws.Range("D36") "18.08"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  -2.88%  "
$ws.Range("E38").Value = "  -0.58%  "
Set-TextValue $ws.Range("D39") "36.83"
$ws.Range("E39").Value = "  +0.94%  "
Set-TextValue $ws.Range("D40") "1.46"
$ws.Range("E40").Value = "  -0.58%  "
Set-TextValue $ws.Range("D41") "0.779"
$ws.Range("E41").Value = "  -1.45%  "
Set-TextValue $ws.Range("D42") "281.52"
$ws.Range("E42").Value = "  +3.50%  "
Set-TextValue $ws.Range("D43") "5.15"
$ws.Range("E43").Value = "  +2.75%  "
Set-TextValue $ws.Range("D44") "3.44"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("E45").Value = "  +1.23%  "
Set-TextValue $ws.Range("D46") "123.35"
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("E47").Value = "  +1.77%  "
Set-TextValue $ws.Range("D48") "17.89"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("E50").Value = "  +0.05%  "
Set-TextValue $ws.Range("D51") "17.09"
$ws.Range("E51").Value = "  +0.11%  "
